$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the CodeSystem version (row 3, column B: "Version" | "0.4.0" -> "0.7.0")
$ws.Range("B3").Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" metadata row entirely (was row 11);
# everything below shifts up by one row.
$ws.Rows.Item(11).Delete()
